$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Formula = '="state"'
$ws.Range("B1").Formula = '="P"'
$ws.Range("C1").Formula = '="V"'
$ws.Range("D1").Formula = '="T"'

$ws.Range("D2").Select()
